# Regenerate merged AHB files
# - rename the "_old" / "_new" suffixed header columns to "_FV2410" / "_FV2504"
# - wrap the data range in an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row labels: A1:J1 "<name>_old" -> "<name>_FV2410"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2410')
}

# K1 ("diff") stays as-is.

# L1:U1 "<name>_new" -> "<name>_FV2504"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2504')
}

# 2) Turn the A1:U57 range into an Excel Table (ListObject) with headers.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the header row (split under row 1).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
